$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.403.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.550.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.74%  '
$ws.Range("E6").Value = '  -2.18%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.99'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.772.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.556.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.344.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("E15").Value = '  -2.76%  '
$ws.Range("E16").Value = '  -2.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0674'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  -3.36%  '
$ws.Range("E30").Value = '  -3.48%  '
$ws.Range("E31").Value = '  -4.95%  '
$ws.Range("E32").Value = '  -1.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.384.42'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("E34").Value = '  -3.67%  '
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("E36").Value = '  -3.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.54%  '
$ws.Range("E38").Value = '  -3.71%  '
$ws.Range("E39").Value = '  -2.89%  '
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.509'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.772'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("E44").Value = '  -2.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.685.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.01%  '
$ws.Range("E48").Value = '  -9.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0101'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.66%  '
